$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "2007年" data row); remaining rows shift up.
$ws.Rows.Item(2).Delete()
